$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 10 (pushes old rows 10-12 down to rows 12-14)
$ws.Rows("10:11").Insert()

# --- Update existing rows 3-9 with revised odds ---
$row3 = @{"F"=1.72; "G"=1.9; "H"=4.3; "I"=5.8; "K"=5.1; "N"=4.2; "O"=1.24; "P"=2.14; "Q"=1.84; "R"=1.45; "S"=2.52; "T"=1.7; "U"=2.12; "W"=2.1}
foreach ($col in $row3.Keys) { $ws.Range($col + "3").Value = $row3[$col] }

$row4 = @{"Q"=1.76; "Z"=40}
foreach ($col in $row4.Keys) { $ws.Range($col + "4").Value = $row4[$col] }

$row5 = @{"L"=1.46; "Z"=980; "AN"=21}
foreach ($col in $row5.Keys) { $ws.Range($col + "5").Value = $row5[$col] }

$row6 = @{"F"=2.38; "G"=2.64; "H"=3.55; "I"=4.3; "J"=2.74; "K"=3.25; "L"=1.61; "M"=1.14; "N"=2.3; "O"=1.58; "T"=2.2; "U"=1.67; "V"=1.3; "W"=1.6; "X"=7.8; "Y"=12; "Z"=29; "AA"=110; "AB"=7.4; "AC"=7.4; "AD"=19.5; "AE"=80; "AF"=15; "AG"=13.5; "AI"=120; "AJ"=44; "AK"=44; "AL"=80; "AM"=280; "AO"=130}
foreach ($col in $row6.Keys) { $ws.Range($col + "6").Value = $row6[$col] }

$row7 = @{"I"=2.76; "O"=1.64; "V"=1.58; "W"=1.28}
foreach ($col in $row7.Keys) { $ws.Range($col + "7").Value = $row7[$col] }

$row8 = @{"F"=2.04; "I"=4.9; "V"=1.26}
foreach ($col in $row8.Keys) { $ws.Range($col + "8").Value = $row8[$col] }

$row9 = @{"F"=1.88; "G"=2.1; "H"=4.8; "I"=6.4; "J"=3.05; "K"=3.85; "L"=1.49; "M"=1.12; "N"=2.44; "O"=1.55; "P"=1.49; "Q"=2.62; "R"=1.17; "S"=4.9; "T"=2.26; "U"=1.65; "V"=1.2; "W"=1.92; "X"=980; "Y"=980; "Z"=980; "AA"=1000; "AB"=7.4; "AC"=980; "AD"=980; "AE"=1000; "AF"=980; "AG"=980; "AH"=980; "AI"=1000; "AJ"=980; "AK"=980; "AL"=80; "AM"=1000; "AN"=980; "AO"=1000}
foreach ($col in $row9.Keys) { $ws.Range($col + "9").Value = $row9[$col] }

# --- New row 10: Cypriot 1st Division ---
$ws.Range("A10").Value = "Cypriot 1st Division"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "2026-01-09"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "14:00:00"
$ws.Range("D10").Value = "A.E.L."
$ws.Range("E10").Value = "Omonia"
$row10 = @{"F"=5.1; "G"=7.2; "H"=1.63; "I"=1.8; "J"=3.9; "K"=4.6; "L"=1.01; "M"=1.01; "N"=2.04; "O"=1.26; "P"=2.04; "Q"=1.78; "R"=1.36; "S"=2.5; "T"=1.01; "U"=1.01; "V"=2.24; "W"=1.16; "X"=1000; "Y"=1000; "Z"=1000; "AA"=1000; "AB"=1000; "AC"=1000; "AD"=1000; "AE"=1000; "AF"=1000; "AG"=1000; "AH"=1000; "AI"=1000; "AJ"=1000; "AK"=1000; "AL"=1000; "AM"=1000; "AN"=1000; "AO"=1000}
foreach ($col in $row10.Keys) { $ws.Range($col + "10").Value = $row10[$col] }

# --- New row 11: Turkish 1 Lig ---
$ws.Range("A11").Value = "Turkish 1 Lig"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2026-01-09"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "14:00:00"
$ws.Range("D11").Value = "Amed Sportif Faaliyetler"
$ws.Range("E11").Value = "Corum Belediyespor"
$row11 = @{"F"=1.79; "G"=2.1; "H"=4.3; "I"=5.8; "J"=3.45; "K"=4.3; "L"=1.32; "M"=1.01; "N"=1.97; "O"=1.28; "P"=1.97; "Q"=1.82; "R"=1.32; "S"=2.74; "T"=1.52; "U"=1.74; "V"=1.2; "W"=1.91; "X"=23; "Y"=980; "Z"=50; "AA"=1000; "AB"=13.5; "AC"=13; "AD"=980; "AE"=90; "AF"=17.5; "AG"=15; "AH"=980; "AI"=90; "AJ"=980; "AK"=980; "AL"=50; "AM"=1000; "AN"=1000; "AO"=1000}
foreach ($col in $row11.Keys) { $ws.Range($col + "11").Value = $row11[$col] }

# --- Updates to shifted rows: 12 (was row 10), 13 (was row 11), 14 (was row 12) ---
$row12 = @{"L"=1.01; "M"=1.04; "N"=5; "O"=1.19; "Q"=1.58; "R"=1.56; "S"=2.34; "T"=1.67; "U"=1.01; "V"=2.92; "W"=1.13; "X"=1000; "Y"=1000; "Z"=14.5; "AA"=18.5; "AB"=44; "AC"=16.5; "AD"=14.5; "AE"=1000; "AF"=95; "AG"=42; "AH"=30; "AI"=44; "AJ"=1000; "AK"=1000; "AL"=1000; "AM"=1000; "AN"=1000; "AO"=1000}
foreach ($col in $row12.Keys) { $ws.Range($col + "12").Value = $row12[$col] }

$row13 = @{"N"=5.1}
foreach ($col in $row13.Keys) { $ws.Range($col + "13").Value = $row13[$col] }

$row14 = @{"F"=3.65; "G"=3.75; "H"=2.5; "I"=2.54; "T"=2.5; "U"=1.63; "Y"=6.4; "AB"=8.8; "AJ"=95}
foreach ($col in $row14.Keys) { $ws.Range($col + "14").Value = $row14[$col] }
